$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range("D2").Value = "27.272.69"
$ws.Range("E2").Value = "  -2.95%  "
$ws.Range("D3").Value = "1.851.14"
$ws.Range("E3").Value = "  -3.65%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "324.87"
$ws.Range("E5").Value = "  -1.84%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "0.4546"
$ws.Range("E7").Value = "  -3.51%  "
$ws.Range("D8").Value = "0.3877"
$ws.Range("E8").Value = "  -4.08%  "
$ws.Range("D9").Value = "48.31"
$ws.Range("E9").Value = "  -9.06%  "
$ws.Range("D10").Value = "0.07906"
$ws.Range("E10").Value = "  -6.19%  "
$ws.Range("D11").Value = "1.012"
$ws.Range("E11").Value = "  -3.22%  "
$ws.Range("D12").Value = "21.33"
$ws.Range("E12").Value = "  -3.81%  "
$ws.Range("D13").Value = "1.859.46"
$ws.Range("E13").Value = "  -3.13%  "
$ws.Range("D14").Value = "5.897"
$ws.Range("E14").Value = "  -3.01%  "
$ws.Range("D15").Value = "7.131"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "0.06596"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").Value = "85.70"
$ws.Range("E18").Value = "  -4.94%  "
$ws.Range("D19").Value = "0.00001025"
$ws.Range("E19").Value = "  -3.88%  "
$ws.Range("D20").Value = "17.11"
$ws.Range("E20").Value = "  -5.16%  "
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "5.483"
$ws.Range("E22").Value = "  -4.37%  "
$ws.Range("D23").Value = "27.284.70"
$ws.Range("E23").Value = "  -2.90%  "
$ws.Range("D24").Value = "10.82"
$ws.Range("E24").Value = "  -4.74%  "
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("D26").Value = "2.081.42"
$ws.Range("E26").Value = "  -2.99%  "
$ws.Range("D27").Value = "154.22"
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").Value = "19.89"
$ws.Range("E29").Value = "  -4.17%  "
$ws.Range("D30").Value = "5.429"
$ws.Range("E30").Value = "  -5.03%  "
$ws.Range("D31").Value = "121.10"
$ws.Range("E31").Value = "  -1.97%  "
$ws.Range("D32").Value = "0.09318"
$ws.Range("E32").Value = "  -2.91%  "
$ws.Range("D33").Value = "0.9400"
$ws.Range("E33").Value = "  -3.55%  "
$ws.Range("D34").Value = "1.444"
$ws.Range("E34").Value = "  -0.38%  "
$ws.Range("D35").Value = "3.588"
$ws.Range("E35").Value = "  -1.45%  "
$ws.Range("D36").Value = "5.246"
$ws.Range("E36").Value = "  -5.34%  "
$ws.Range("D37").Value = "0.06027"
$ws.Range("E37").Value = "  -2.03%  "
$ws.Range("D38").Value = "0.02222"
$ws.Range("E38").Value = "  -3.71%  "
$ws.Range("D39").Value = "1.204"
$ws.Range("E39").Value = "  -2.36%  "
$ws.Range("D40").Value = "8.054"
$ws.Range("E40").Value = "  -9.92%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").Value = "0.5913"
$ws.Range("E42").Value = "  -3.77%  "
$ws.Range("D43").Value = "0.1882"
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("D44").Value = "10.14"
$ws.Range("E44").Value = "  -8.21%  "
$ws.Range("D45").Value = "1.269"
$ws.Range("E45").Value = "  -2.61%  "
$ws.Range("D46").Value = "0.5599"
$ws.Range("E46").Value = "  -4.54%  "
$ws.Range("D47").Value = "12.10"
$ws.Range("E47").Value = "  -5.46%  "
$ws.Range("D48").Value = "3.377"
$ws.Range("E48").Value = "  -2.88%  "
$ws.Range("D49").Value = "1.907"
$ws.Range("E49").Value = "  -6.08%  "
$ws.Range("D50").Value = "0.06729"
$ws.Range("E50").Value = "  -1.53%  "
$ws.Range("D51").Value = "107.69"
$ws.Range("E51").Value = "  -2.13%  "

$priceCol.Style = "Normal"
